# Updates cryptos list values (price + 1h volume change) per upstream refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.134.37"
$ws.Range("E2").Value = "  -3.88%  "

$ws.Range("D3").Value = "3.500.84"
$ws.Range("E3").Value = "  -5.84%  "

$ws.Range("D5").Value = "'581.26"
$ws.Range("E5").Value = "  -1.40%  "

$ws.Range("D6").Value = "'172.93"
$ws.Range("E6").Value = "  -4.41%  "

$ws.Range("D7").Value = "'0.618"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "3.491.52"
$ws.Range("E8").Value = "  -5.92%  "

$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("D10").Value = "'0.188"
$ws.Range("E10").Value = "  -7.56%  "

$ws.Range("D11").Value = "'6.69"
$ws.Range("E11").Value = "  +4.29%  "

$ws.Range("D12").Value = "'0.593"
$ws.Range("E12").Value = "  -3.43%  "

$ws.Range("D13").Value = "'46.69"
$ws.Range("E13").Value = "  -6.71%  "

$ws.Range("D14").Value = "'0.0000275"
$ws.Range("E14").Value = "  -4.58%  "

$ws.Range("D15").Value = "'673.60"
$ws.Range("E15").Value = "  -1.44%  "

$ws.Range("D16").Value = "4.061.68"
$ws.Range("E16").Value = "  -5.90%  "

$ws.Range("D17").Value = "'8.68"
$ws.Range("E17").Value = "  -4.20%  "

$ws.Range("D18").Value = "69.139.58"
$ws.Range("E18").Value = "  -3.92%  "

$ws.Range("D19").Value = "3.504.02"
$ws.Range("E19").Value = "  -5.68%  "

$ws.Range("E20").Value = "  -1.53%  "

$ws.Range("D21").Value = "'17.38"
$ws.Range("E21").Value = "  -4.44%  "

$ws.Range("D22").Value = "'11.15"
$ws.Range("E22").Value = "  -4.74%  "

$ws.Range("D23").Value = "'0.899"
$ws.Range("E23").Value = "  -5.19%  "

$ws.Range("D24").Value = "'16.09"
$ws.Range("E24").Value = "  -9.78%  "

$ws.Range("D25").Value = "'97.45"
$ws.Range("E25").Value = "  -6.19%  "

$ws.Range("D26").Value = "'3.86"
$ws.Range("E26").Value = "  -4.94%  "

$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.05%  "

$ws.Range("D28").Value = "'2.64"
$ws.Range("E28").Value = "  -7.15%  "

$ws.Range("D29").Value = "'9.38"
$ws.Range("E29").Value = "  -9.50%  "

$ws.Range("D30").Value = "'32.85"
$ws.Range("E30").Value = "  -7.89%  "

$ws.Range("D31").Value = "'8.67"
$ws.Range("E31").Value = "  -6.87%  "

$ws.Range("D32").Value = "'3.17"
$ws.Range("E32").Value = "  -8.63%  "

$ws.Range("E33").Value = "  -6.29%  "

$ws.Range("D34").Value = "'7.24"
$ws.Range("E34").Value = "  -1.94%  "

$ws.Range("D35").Value = "'598.08"
$ws.Range("E35").Value = "  +6.22%  "

$ws.Range("D36").Value = "'3.58"
$ws.Range("E36").Value = "  -15.25%  "

$ws.Range("D37").Value = "'10.84"
$ws.Range("E37").Value = "  -4.14%  "

$ws.Range("E38").Value = "  -5.38%  "

$ws.Range("D39").Value = "'57.04"
$ws.Range("E39").Value = "  -4.48%  "

$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  +0.18%  "

$ws.Range("D41").Value = "'0.0436"
$ws.Range("E41").Value = "  -6.23%  "

$ws.Range("D42").Value = "'0.334"
$ws.Range("E42").Value = "  -5.50%  "

$ws.Range("D43").Value = "3.413.02"
$ws.Range("E43").Value = "  -9.61%  "

$ws.Range("D44").Value = "'0.135"
$ws.Range("E44").Value = "  -6.43%  "

$ws.Range("D45").Value = "'33.22"
$ws.Range("E45").Value = "  -7.40%  "

$ws.Range("D46").Value = "0.0₃0705"
$ws.Range("E46").Value = "  -9.63%  "

$ws.Range("D47").Value = "'2.88"
$ws.Range("E47").Value = "  -0.78%  "

$ws.Range("D48").Value = "'2.59"
$ws.Range("E48").Value = "  -7.90%  "

$ws.Range("E49").Value = "  -0.95%  "

$ws.Range("D50").Value = "'5.75"
$ws.Range("E50").Value = "  +17.00%  "

$ws.Range("D51").Value = "'133.16"
$ws.Range("E51").Value = "  -2.05%  "
